$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rerunning the model with correctly initialised weights: the "weight"
# columns now hold clean +/-1 (or 0/1) arrays and every "bias" column is a
# flat 0.0, reflecting the freshly (correctly) initialised parameters.
# (row 2 holds the epoch-0 snapshot of the run)

$ws.Range("C2").Value = "[0. 1.]"
$ws.Range("E2").Value = "[1. 0.]"
$ws.Range("G2").Value = "[1. 1.]"
$ws.Range("I2").Value = "[1. 1.]"
$ws.Range("K2").Value = "[-1.  1.]"
$ws.Range("M2").Value = "[1.]"
$ws.Range("O2").Value = "[ 1. -1.]"
$ws.Range("Q2").Value = "[1. 1.]"
$ws.Range("S2").Value = "[[ 1.  1.]" + [char]10 + " [-1. -1.]]"
$ws.Range("T2").Value = "[0. 1.]"

# The bias columns need to hold the literal text "0.0" (not the number 0),
# so force a text format before assigning, then clear the style again so no
# stray formatting is left behind on the cell.
$biasCells = @("D2","F2","H2","J2","L2","N2","P2","R2")
foreach ($cell in $biasCells) {
    $ws.Range($cell).NumberFormat = "@"
}
foreach ($cell in $biasCells) {
    $ws.Range($cell).Value = "0.0"
}
foreach ($cell in $biasCells) {
    $ws.Range($cell).Style = "Normal"
}

# "all losses" goes back to 0 for the fresh run
$ws.Range("U2").Value = 0

# Undo the stray autofit row-height bump caused by the embedded newline in S2
$ws.Rows(2).AutoFit()
